# Apply crime data update for 2024-11-29
# Updates K-column (year 2024) totals across Citywide Totals, By Neighborhood,
# and per-neighborhood sheets to reflect newly added incident records.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7323
$ws.Range("K3").Value = 7596
$ws.Range("K4").Value = 1585
$ws.Range("K6").Value = 8409
$ws.Range("K7").Value = 25447

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 128
$ws.Range("K7").Value = 312

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 459
$ws.Range("K3").Value = 499
$ws.Range("K6").Value = 555
$ws.Range("K7").Value = 1654

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 269
$ws.Range("K7").Value = 1073

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 424

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 243
$ws.Range("K3").Value = 280
$ws.Range("K7").Value = 848

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 196
$ws.Range("K7").Value = 603

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 770
$ws.Range("K8").Value = 1654
$ws.Range("K9").Value = 118
$ws.Range("K12").Value = 46
$ws.Range("K13").Value = 37
$ws.Range("K19").Value = 735
$ws.Range("K27").Value = 249
$ws.Range("K28").Value = 11
$ws.Range("K29").Value = 1399
$ws.Range("K30").Value = 103
$ws.Range("K33").Value = 1073
$ws.Range("K36").Value = 324
$ws.Range("K37").Value = 848
$ws.Range("K42").Value = 937
$ws.Range("K47").Value = 174
$ws.Range("K51").Value = 322
$ws.Range("K52").Value = 666
$ws.Range("K53").Value = 312
$ws.Range("K55").Value = 275
$ws.Range("K60").Value = 147
$ws.Range("K63").Value = 70
$ws.Range("K65").Value = 603
$ws.Range("K67").Value = 992
$ws.Range("K71").Value = 78
$ws.Range("K78").Value = 310
$ws.Range("K79").Value = 625
$ws.Range("K80").Value = 95
$ws.Range("K82").Value = 33
$ws.Range("K84").Value = 201
$ws.Range("K85").Value = 1161
$ws.Range("K89").Value = 383
$ws.Range("K95").Value = 424
$ws.Range("K101").Value = 25447

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K4").Value = 54
$ws.Range("K7").Value = 992

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 66
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 398
$ws.Range("K3").Value = 497
$ws.Range("K6").Value = 406
$ws.Range("K7").Value = 1399

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 215
$ws.Range("K7").Value = 735

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K4").Value = 43
$ws.Range("K7").Value = 937

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("K2").Value = 6
$ws.Range("K5").Value = 12
$ws.Range("K6").Value = 37

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 81
$ws.Range("K7").Value = 310

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 275

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 205
$ws.Range("K3").Value = 200
$ws.Range("K4").Value = 41
$ws.Range("K7").Value = 625

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 102
$ws.Range("K7").Value = 324

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K6").Value = 214
$ws.Range("K7").Value = 770

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 174

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 107
$ws.Range("K4").Value = 42
$ws.Range("K6").Value = 115
$ws.Range("K7").Value = 383

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 249

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K4").Value = 36
$ws.Range("K7").Value = 322

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 379
$ws.Range("K7").Value = 1161

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K3").Value = 7
$ws.Range("K6").Value = 33

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 178
$ws.Range("K3").Value = 184
$ws.Range("K4").Value = 39
$ws.Range("K6").Value = 243
$ws.Range("K7").Value = 666

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 11
